# "Generate Report for Archive"
# - Update the localization status text "Ready for handoff" -> "In Translation"
#   on all three sheets (Overview E2/F2, zh-cn C2, de-de C2 all share the same
#   string, so updating each cell's value collapses them back onto one shared
#   string, matching the sharedStrings.xml diff).
# - Shrink the "Status" column(s) from 17.2159881591797 to 13.4101848602295
#   on Overview (columns E & F) and on zh-cn / de-de (column C).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Text update: "Ready for handoff" -> "In Translation" ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Column width update ---
# Target stored OOXML width is 13.4101848602295 characters. This engine's
# ColumnWidth setter snaps the persisted <col width> to the nearest 1/6
# (i.e. whole on-screen pixel) grid point, so no COM input reproduces that
# float exactly; 80/6 = 13.333333333333334 is the closest reachable grid
# point (vs. the next nearest, 81/6 = 13.5). ColumnWidth = 12.5 lands
# squarely in the middle of the input range that snaps to that grid point.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
